# Scheduled market-data refresh: updates live price columns (H-N) across
# the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# row 17
$ws.Range("H17").Value = 747.7846
$ws.Range("J17").Value = 747.7846
$ws.Range("L17").Value = 2243.3538
$ws.Range("N17").Value = -2579.3538
# row 76
$ws.Range("H76").Value = 3649.52
$ws.Range("I76").Value = 3325
$ws.Range("J76").Value = 3802.2354
$ws.Range("K76").Value = 3325
$ws.Range("L76").Value = 3802.2354
$ws.Range("M76").Value = -3010
$ws.Range("N76").Value = -4432.2354
# row 79
$ws.Range("H79").Value = 3649.52
$ws.Range("I79").Value = 3325
$ws.Range("J79").Value = 3802.2354
$ws.Range("K79").Value = 3325
$ws.Range("L79").Value = 3802.2354
$ws.Range("M79").Value = -2233
$ws.Range("N79").Value = -5986.2354
# row 132
$ws.Range("H132").Value = 664060.4
$ws.Range("I132").Value = 1798.9108
$ws.Range("J132").Value = 2724429.5
$ws.Range("K132").Value = 5396.732400000001
$ws.Range("L132").Value = 8173288.5
$ws.Range("M132").Value = -2866.732400000001
$ws.Range("N132").Value = -8178348.5
# row 138
$ws.Range("H138").Value = 1986185.5
$ws.Range("I138").Value = 1282.24
$ws.Range("J138").Value = 4905161
$ws.Range("K138").Value = 3846.72
$ws.Range("L138").Value = 14715483
$ws.Range("M138").Value = 1293.28
$ws.Range("N138").Value = -14725763

$ws = $wb.Worksheets.Item("ARM")
# row 32
$ws.Range("H32").Value = 1985.97
$ws.Range("I32").Value = 1725.0532
$ws.Range("J32").Value = 6073.6665
$ws.Range("K32").Value = 1725.0532
$ws.Range("L32").Value = 6073.6665
$ws.Range("M32").Value = -1438.0532
$ws.Range("N32").Value = -6647.6665
# row 37
$ws.Range("H37").Value = 0
$ws.Range("I37").Value = 0
$ws.Range("K37").Value = 0
$ws.Range("M37").ClearContents()
# row 44
$ws.Range("H44").Value = 20000
$ws.Range("J44").Value = 20000
$ws.Range("L44").Value = 20000
$ws.Range("N44").Value = -20976
# row 55
$ws.Range("H55").Value = 20000
$ws.Range("J55").Value = 20000
$ws.Range("L55").Value = 20000
$ws.Range("N55").Value = -20630
# row 61
$ws.Range("H61").Value = 20876520
$ws.Range("I61").Value = 24415818
$ws.Range("J61").Value = 146346.86
$ws.Range("K61").Value = 24415818
$ws.Range("L61").Value = 146346.86
$ws.Range("M61").Value = -24415606
$ws.Range("N61").Value = -146770.86
# row 74
$ws.Range("H74").Value = 11179464
$ws.Range("I74").Value = 16718143
$ws.Range("J74").Value = 102106.3
$ws.Range("K74").Value = 16718143
$ws.Range("L74").Value = 102106.3
$ws.Range("M74").Value = -16717269
$ws.Range("N74").Value = -103854.3
# row 77
$ws.Range("H77").Value = 11179464
$ws.Range("I77").Value = 16718143
$ws.Range("J77").Value = 102106.3
$ws.Range("K77").Value = 83590715
$ws.Range("L77").Value = 510531.5
$ws.Range("M77").Value = -83586347
$ws.Range("N77").Value = -519267.5
# row 80
$ws.Range("H80").Value = 34985
$ws.Range("I80").Value = 35000
$ws.Range("J80").Value = 34980
$ws.Range("K80").Value = 35000
$ws.Range("L80").Value = 34980
$ws.Range("M80").Value = -34002
$ws.Range("N80").Value = -36976
# row 83
$ws.Range("H83").Value = 34985
$ws.Range("I83").Value = 35000
$ws.Range("J83").Value = 34980
$ws.Range("K83").Value = 105000
$ws.Range("L83").Value = 104940
$ws.Range("M83").Value = -100008
$ws.Range("N83").Value = -114924
# row 132
$ws.Range("H132").Value = 35705.918
$ws.Range("I132").Value = 24441.785
$ws.Range("J132").Value = 61988.89
$ws.Range("K132").Value = 73325.355
$ws.Range("L132").Value = 185966.67
$ws.Range("M132").Value = -70795.355
$ws.Range("N132").Value = -191026.67
# row 136
$ws.Range("H136").Value = 20876520
$ws.Range("I136").Value = 24415818
$ws.Range("J136").Value = 146346.86
$ws.Range("K136").Value = 73247454
$ws.Range("L136").Value = 439040.58
$ws.Range("M136").Value = -73244904
$ws.Range("N136").Value = -444140.58

$ws = $wb.Worksheets.Item("BSM")
# row 105
$ws.Range("H105").Value = 15627113
$ws.Range("I105").Value = 27779752
$ws.Range("J105").Value = 2292.2144
$ws.Range("K105").Value = 27779752
$ws.Range("L105").Value = 2292.2144
$ws.Range("M105").Value = -27778005
$ws.Range("N105").Value = -5786.2144
# row 134
$ws.Range("H134").Value = 1787.5745
$ws.Range("I134").Value = 983.7742
$ws.Range("J134").Value = 3344.9375
$ws.Range("K134").Value = 2951.3226
$ws.Range("L134").Value = 10034.8125
$ws.Range("M134").Value = -416.3226
$ws.Range("N134").Value = -15104.8125

$ws = $wb.Worksheets.Item("CRP")
# row 31
$ws.Range("H31").Value = 3153.1
$ws.Range("I31").Value = 1627.8462
$ws.Range("J31").Value = 5985.7144
$ws.Range("K31").Value = 1627.8462
$ws.Range("L31").Value = 5985.7144
$ws.Range("M31").Value = -1332.8462
$ws.Range("N31").Value = -6575.7144
# row 34
$ws.Range("H34").Value = 3153.1
$ws.Range("I34").Value = 1627.8462
$ws.Range("J34").Value = 5985.7144
$ws.Range("K34").Value = 1627.8462
$ws.Range("L34").Value = 5985.7144
$ws.Range("M34").Value = -1425.8462
$ws.Range("N34").Value = -6389.7144
# row 58
$ws.Range("H58").Value = 23811478
$ws.Range("I58").Value = 31251838
$ws.Range("J58").Value = 2330.8
$ws.Range("K58").Value = 31251838
$ws.Range("L58").Value = 2330.8
$ws.Range("M58").Value = -31251635
$ws.Range("N58").Value = -2736.8
# row 132
$ws.Range("H132").Value = 18148.268
$ws.Range("I132").Value = 1302.6222
$ws.Range("J132").Value = 68685.2
$ws.Range("K132").Value = 3907.8666
$ws.Range("L132").Value = 206055.6
$ws.Range("M132").Value = -1377.8666
$ws.Range("N132").Value = -211115.6
# row 136
$ws.Range("H136").Value = 23811478
$ws.Range("I136").Value = 31251838
$ws.Range("J136").Value = 2330.8
$ws.Range("K136").Value = 93755514
$ws.Range("L136").Value = 6992.400000000001
$ws.Range("M136").Value = -93752964
$ws.Range("N136").Value = -12092.4

$ws = $wb.Worksheets.Item("CUL")
# row 98
$ws.Range("H98").Value = 646.375
$ws.Range("J98").Value = 710.1429000000001
$ws.Range("L98").Value = 2130.4287
$ws.Range("N98").Value = -5126.4287
# row 100
$ws.Range("H100").Value = 2773.913
$ws.Range("J100").Value = 2773.913
$ws.Range("L100").Value = 8321.739
$ws.Range("N100").Value = -9943.739
# row 104
$ws.Range("H104").Value = 4000
$ws.Range("J104").Value = 4000
$ws.Range("L104").Value = 12000
$ws.Range("N104").Value = -17242
# row 131
$ws.Range("H131").Value = 1292.9259
$ws.Range("I131").Value = 380
$ws.Range("J131").Value = 1500.409
$ws.Range("K131").Value = 1140
$ws.Range("L131").Value = 4501.227000000001
$ws.Range("M131").Value = 3900
$ws.Range("N131").Value = -14581.227

$ws = $wb.Worksheets.Item("GSM")
# row 43
$ws.Range("H43").Value = 5458
$ws.Range("I43").Value = 650
$ws.Range("J43").Value = 8663.333000000001
$ws.Range("K43").Value = 650
$ws.Range("L43").Value = 8663.333000000001
$ws.Range("M43").Value = -499
$ws.Range("N43").Value = -8965.333000000001
# row 80
$ws.Range("H80").Value = 3579.625
$ws.Range("I80").Value = 3048.9
$ws.Range("K80").Value = 3048.9
$ws.Range("M80").Value = -2050.9
# row 83
$ws.Range("H83").Value = 3579.625
$ws.Range("I83").Value = 3048.9
$ws.Range("K83").Value = 15244.5
$ws.Range("M83").Value = -10252.5
# row 132
$ws.Range("H132").Value = 30759.795
$ws.Range("I132").Value = 20347.576
$ws.Range("J132").Value = 64599.5
$ws.Range("K132").Value = 61042.728
$ws.Range("L132").Value = 193798.5
$ws.Range("M132").Value = -58512.728
$ws.Range("N132").Value = -198858.5

$ws = $wb.Worksheets.Item("LTW")
# row 100
$ws.Range("H100").Value = 1410.2307
$ws.Range("I100").Value = 1181.4445
$ws.Range("K100").Value = 1181.4445
$ws.Range("M100").Value = -640.4445000000001

$ws = $wb.Worksheets.Item("WVR")
# row 122
$ws.Range("H122").Value = 1626.2778
$ws.Range("I122").Value = 1179.4736
$ws.Range("J122").Value = 2687.4375
$ws.Range("K122").Value = 3538.4208
$ws.Range("L122").Value = 8062.3125
$ws.Range("M122").Value = -1088.4208
$ws.Range("N122").Value = -12962.3125
# row 126
$ws.Range("H126").Value = 1149.75
$ws.Range("I126").Value = 1166.3889
$ws.Range("K126").Value = 3499.1667
$ws.Range("M126").Value = -1029.1667
# row 132
$ws.Range("H132").Value = 45647.42
$ws.Range("I132").Value = 28560.723
$ws.Range("J132").Value = 113994.22
$ws.Range("K132").Value = 85682.16900000001
$ws.Range("L132").Value = 341982.66
$ws.Range("M132").Value = -83152.16900000001
$ws.Range("N132").Value = -347042.66
# row 136
$ws.Range("H136").Value = 50002.668
$ws.Range("I136").Value = 38109.703
$ws.Range("J136").Value = 71410
$ws.Range("K136").Value = 114329.109
$ws.Range("L136").Value = 214230
$ws.Range("M136").Value = -111779.109
$ws.Range("N136").Value = -219330
